$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new status labels below the existing list (rows 66-70)
$ws.Range("A67").Value = "Pesanan Dikomplain"
$ws.Range("A68").Value = "Pesanan Selesai"
$ws.Range("A69").Value = "Pesanan Tiba"
$ws.Range("A70").Value = "Resi Diubah"
$ws.Range("A66").Value = "Pesanan Dikirim"

# Match the author's final selection/scroll position (A66 becomes the active cell)
$ws.Range("A66").Select()
